$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.797.77'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +8.70%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.951.44'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +7.24%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9996'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.27%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '342.31'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +2.91%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9998'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.30%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4764'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +4.45%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4151'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +9.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '48.43'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +5.68%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08276'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +5.66%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.042'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +9.20%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.71'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +8.42%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.214'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +6.87%  '

$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.936.43'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +6.04%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.435'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +5.52%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '92.33'

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.001'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.24%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001067'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +4.90%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06676'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.78%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.10'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +6.26%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.11%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '29.738.55'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +8.59%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.604'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +6.24%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.29'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +4.67%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.283'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.37%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.170.97'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +6.17%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '161.04'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +1.96%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.23'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +5.07%  '

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +7.90%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.704'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +8.75%  '

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +4.41%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.026'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +10.51%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09632'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +3.53%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.485'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +13.12%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.682'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +3.28%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.524'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +6.31%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06326'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +7.45%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02334'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +6.99%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.581'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +6.43%  '

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +5.34%  '

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +6.89%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '10.76'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +8.86%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1902'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +5.12%  '

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.22%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.274'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.23%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.392'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +33.01%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.62'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +7.08%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5736'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +6.73%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.005'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +7.41%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07400'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +12.77%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '114.23'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +3.86%  '
